$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 277.24
$ws.Cells.Item(15, 9).Value = 277.24
$ws.Cells.Item(15, 11).Value = 831.72
$ws.Cells.Item(15, 13).Value = -662.72

$ws.Cells.Item(33, 8).Value = 197.5
$ws.Cells.Item(33, 9).Value = 160
$ws.Cells.Item(33, 11).Value = 160
$ws.Cells.Item(33, 13).Value = 69

$ws.Cells.Item(87, 8).Value = 13731.875
$ws.Cells.Item(87, 10).Value = 13731.875
$ws.Cells.Item(87, 12).Value = 13731.875
$ws.Cells.Item(87, 14).Value = -16227.875

$ws.Cells.Item(90, 8).Value = 13731.875
$ws.Cells.Item(90, 10).Value = 13731.875
$ws.Cells.Item(90, 12).Value = 41195.625
$ws.Cells.Item(90, 14).Value = -53675.625

$ws.Cells.Item(92, 8).Value = 3380.0527
$ws.Cells.Item(92, 9).Value = 3745.0588
$ws.Cells.Item(92, 11).Value = 3745.0588
$ws.Cells.Item(92, 13).Value = -2497.0588

$ws.Cells.Item(132, 8).Value = 2718435.8
$ws.Cells.Item(132, 9).Value = 3206202.2
$ws.Cells.Item(132, 10).Value = 878.5714
$ws.Cells.Item(132, 11).Value = 9618606.600000001
$ws.Cells.Item(132, 12).Value = 2635.7142
$ws.Cells.Item(132, 13).Value = -9616076.600000001
$ws.Cells.Item(132, 14).Value = -7695.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1120.871
$ws.Cells.Item(2, 9).Value = 844.1923
$ws.Cells.Item(2, 10).Value = 2559.6
$ws.Cells.Item(2, 11).Value = 844.1923
$ws.Cells.Item(2, 12).Value = 2559.6
$ws.Cells.Item(2, 13).Value = -731.1923
$ws.Cells.Item(2, 14).Value = -2785.6

$ws.Cells.Item(32, 8).Value = 1316.04
$ws.Cells.Item(32, 9).Value = 1326.3131
$ws.Cells.Item(32, 10).Value = 299
$ws.Cells.Item(32, 11).Value = 1326.3131
$ws.Cells.Item(32, 12).Value = 299
$ws.Cells.Item(32, 13).Value = -1039.3131
$ws.Cells.Item(32, 14).Value = -873

$ws.Cells.Item(45, 8).Value = 1094.1277
$ws.Cells.Item(45, 9).Value = 1124.079
$ws.Cells.Item(45, 10).Value = 967.6667
$ws.Cells.Item(45, 11).Value = 1124.079
$ws.Cells.Item(45, 12).Value = 967.6667
$ws.Cells.Item(45, 13).Value = -747.079
$ws.Cells.Item(45, 14).Value = -1721.6667

$ws.Cells.Item(104, 8).Value = 49000
$ws.Cells.Item(104, 10).Value = 49000
$ws.Cells.Item(104, 12).Value = 49000
$ws.Cells.Item(104, 14).Value = -55988

$ws.Cells.Item(116, 8).Value = 1120.871
$ws.Cells.Item(116, 9).Value = 844.1923
$ws.Cells.Item(116, 10).Value = 2559.6
$ws.Cells.Item(116, 11).Value = 844.1923
$ws.Cells.Item(116, 12).Value = 2559.6
$ws.Cells.Item(116, 13).Value = 1449.8077
$ws.Cells.Item(116, 14).Value = -7147.6

$ws.Cells.Item(122, 8).Value = 1940.2
$ws.Cells.Item(122, 9).Value = 1437.75
$ws.Cells.Item(122, 10).Value = 3950
$ws.Cells.Item(122, 11).Value = 4313.25
$ws.Cells.Item(122, 12).Value = 11850
$ws.Cells.Item(122, 13).Value = -1863.25
$ws.Cells.Item(122, 14).Value = -16750

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1120.871
$ws.Cells.Item(3, 9).Value = 844.1923
$ws.Cells.Item(3, 10).Value = 2559.6
$ws.Cells.Item(3, 11).Value = 844.1923
$ws.Cells.Item(3, 12).Value = 2559.6
$ws.Cells.Item(3, 13).Value = -730.1923
$ws.Cells.Item(3, 14).Value = -2787.6

$ws.Cells.Item(94, 8).Value = 942.25
$ws.Cells.Item(94, 9).Value = 742.73914
$ws.Cells.Item(94, 10).Value = 1860
$ws.Cells.Item(94, 11).Value = 742.73914
$ws.Cells.Item(94, 12).Value = 1860
$ws.Cells.Item(94, 13).Value = -291.73914
$ws.Cells.Item(94, 14).Value = -2762

$ws.Cells.Item(105, 8).Value = 2045.7
$ws.Cells.Item(105, 9).Value = 1995.2222
$ws.Cells.Item(105, 10).Value = 2500
$ws.Cells.Item(105, 11).Value = 1995.2222
$ws.Cells.Item(105, 12).Value = 2500
$ws.Cells.Item(105, 13).Value = -248.2221999999999
$ws.Cells.Item(105, 14).Value = -5994

$ws.Cells.Item(107, 8).Value = 1311.4286
$ws.Cells.Item(107, 9).Value = 1247.5
$ws.Cells.Item(107, 10).Value = 1396.6666
$ws.Cells.Item(107, 11).Value = 1247.5
$ws.Cells.Item(107, 12).Value = 1396.6666
$ws.Cells.Item(107, 13).Value = 672.5
$ws.Cells.Item(107, 14).Value = -5236.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 268273.72
$ws.Cells.Item(94, 9).Value = 334349.66
$ws.Cells.Item(94, 10).Value = 224223.11
$ws.Cells.Item(94, 11).Value = 334349.66
$ws.Cells.Item(94, 12).Value = 224223.11
$ws.Cells.Item(94, 13).Value = -333898.66
$ws.Cells.Item(94, 14).Value = -225125.11

$ws.Cells.Item(107, 8).Value = 1737.4706
$ws.Cells.Item(107, 9).Value = 2270.0908
$ws.Cells.Item(107, 10).Value = 761
$ws.Cells.Item(107, 11).Value = 2270.0908
$ws.Cells.Item(107, 12).Value = 761
$ws.Cells.Item(107, 13).Value = -350.0907999999999
$ws.Cells.Item(107, 14).Value = -4601

$ws.Cells.Item(125, 8).Value = 22208.666
$ws.Cells.Item(125, 10).Value = 22208.666
$ws.Cells.Item(125, 12).Value = 22208.666
$ws.Cells.Item(125, 14).Value = -27128.666

$ws.Cells.Item(132, 8).Value = 1431.7778
$ws.Cells.Item(132, 9).Value = 1025.76
$ws.Cells.Item(132, 10).Value = 2354.5454
$ws.Cells.Item(132, 11).Value = 3077.28
$ws.Cells.Item(132, 12).Value = 7063.6362
$ws.Cells.Item(132, 13).Value = -547.2799999999997
$ws.Cells.Item(132, 14).Value = -12123.6362

$ws.Cells.Item(134, 8).Value = 1749.4517
$ws.Cells.Item(134, 9).Value = 1754.25
$ws.Cells.Item(134, 10).Value = 1704.6666
$ws.Cells.Item(134, 11).Value = 5262.75
$ws.Cells.Item(134, 12).Value = 5113.9998
$ws.Cells.Item(134, 13).Value = -2727.75
$ws.Cells.Item(134, 14).Value = -10183.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 1128.7142
$ws.Cells.Item(39, 9).Value = 800
$ws.Cells.Item(39, 10).Value = 1183.5
$ws.Cells.Item(39, 11).Value = 2400
$ws.Cells.Item(39, 12).Value = 3550.5
$ws.Cells.Item(39, 13).Value = -2106
$ws.Cells.Item(39, 14).Value = -4138.5

$ws.Cells.Item(119, 8).Value = 2490.2964
$ws.Cells.Item(119, 9).Value = 1220.4286
$ws.Cells.Item(119, 10).Value = 3857.8462
$ws.Cells.Item(119, 11).Value = 3661.2858
$ws.Cells.Item(119, 12).Value = 11573.5386
$ws.Cells.Item(119, 13).Value = 1176.7142
$ws.Cells.Item(119, 14).Value = -21249.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1276.5
$ws.Cells.Item(97, 9).Value = 949.63635
$ws.Cells.Item(97, 10).Value = 2475
$ws.Cells.Item(97, 11).Value = 949.63635
$ws.Cells.Item(97, 12).Value = 2475
$ws.Cells.Item(97, 13).Value = -453.63635
$ws.Cells.Item(97, 14).Value = -3467

$ws.Cells.Item(107, 8).Value = 5263942.5
$ws.Cells.Item(107, 9).Value = 643.625
$ws.Cells.Item(107, 10).Value = 26317138
$ws.Cells.Item(107, 11).Value = 643.625
$ws.Cells.Item(107, 12).Value = 26317138
$ws.Cells.Item(107, 13).Value = 1276.375
$ws.Cells.Item(107, 14).Value = -26320978

$ws.Cells.Item(113, 8).Value = 2002
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 2002
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 2002
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -6342

$ws.Cells.Item(122, 8).Value = 2864.2727
$ws.Cells.Item(122, 9).Value = 1643.8572
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 4931.571599999999
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -2481.571599999999
$ws.Cells.Item(122, 14).Value = -19900

$ws.Cells.Item(132, 8).Value = 1759
$ws.Cells.Item(132, 9).Value = 1382.4
$ws.Cells.Item(132, 10).Value = 3223.5557
$ws.Cells.Item(132, 11).Value = 4147.200000000001
$ws.Cells.Item(132, 12).Value = 9670.667099999999
$ws.Cells.Item(132, 13).Value = -1617.200000000001
$ws.Cells.Item(132, 14).Value = -14730.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2662.8125
$ws.Cells.Item(40, 9).Value = 2339.4443
$ws.Cells.Item(40, 10).Value = 3078.5715
$ws.Cells.Item(40, 11).Value = 2339.4443
$ws.Cells.Item(40, 12).Value = 3078.5715
$ws.Cells.Item(40, 13).Value = -2203.4443
$ws.Cells.Item(40, 14).Value = -3350.5715

$ws.Cells.Item(122, 8).Value = 3726.5789
$ws.Cells.Item(122, 9).Value = 4177.778
$ws.Cells.Item(122, 11).Value = 12533.334
$ws.Cells.Item(122, 13).Value = -10083.334

$ws.Cells.Item(136, 8).Value = 3472.6
$ws.Cells.Item(136, 9).Value = 3061.3333
$ws.Cells.Item(136, 10).Value = 3809.0908
$ws.Cells.Item(136, 11).Value = 9183.999899999999
$ws.Cells.Item(136, 12).Value = 11427.2724
$ws.Cells.Item(136, 13).Value = -6633.999899999999
$ws.Cells.Item(136, 14).Value = -16527.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 15666.667
$ws.Cells.Item(43, 9).Value = 13500
$ws.Cells.Item(43, 10).Value = 20000
$ws.Cells.Item(43, 11).Value = 13500
$ws.Cells.Item(43, 12).Value = 20000
$ws.Cells.Item(43, 13).Value = -13351
$ws.Cells.Item(43, 14).Value = -20298

$ws.Cells.Item(132, 8).Value = 1391.3103
$ws.Cells.Item(132, 9).Value = 1072.1428
$ws.Cells.Item(132, 10).Value = 1689.2
$ws.Cells.Item(132, 11).Value = 3216.4284
$ws.Cells.Item(132, 12).Value = 5067.6
$ws.Cells.Item(132, 13).Value = -686.4284000000002
$ws.Cells.Item(132, 14).Value = -10127.6
